$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe
# (Excel will store these as text since they don't parse as numbers)
$ws.Range("D2").Value = '68.598.90'
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = '3.841.85'
$ws.Range("E3").Value = '  -2.34%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("E5").Value = '  +4.52%  '
$ws.Range("E6").Value = '  -4.98%  '
$ws.Range("E7").Value = '  -3.04%  '
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("E9").Value = '  -3.29%  '
$ws.Range("E10").Value = '  -4.94%  '
$ws.Range("E11").Value = '  -8.51%  '
$ws.Range("E12").Value = '  -4.42%  '
$ws.Range("E13").Value = '  -2.02%  '
$ws.Range("D14").Value = '4.460.76'
$ws.Range("E14").Value = '  -2.41%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.890.38'
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("E16").Value = '  +6.06%  '
$ws.Range("E17").Value = '  -3.19%  '
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").Value = '68.538.37'
$ws.Range("E20").Value = '  -1.02%  '
$ws.Range("E21").Value = '  -5.80%  '
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("E24").Value = '  -4.23%  '
$ws.Range("E25").Value = '  -2.66%  '
$ws.Range("E26").Value = '  +4.52%  '
$ws.Range("E27").Value = '  -7.13%  '
$ws.Range("E28").Value = '  -5.11%  '
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("E30").Value = '  -3.82%  '
$ws.Range("E31").Value = '  +13.95%  '
$ws.Range("E32").Value = '  -1.77%  '
$ws.Range("E33").Value = '  -5.55%  '
$ws.Range("E34").Value = '  +7.03%  '
$ws.Range("E35").Value = '  -5.60%  '
$ws.Range("D36").Value = '0.0₃0839'
$ws.Range("E36").Value = '  -6.88%  '
$ws.Range("E37").Value = '  -3.96%  '
$ws.Range("B38").Value = 'ThetaToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("E38").Value = '  +10.00%  '
$ws.Range("E39").Value = '  -3.38%  '
$ws.Range("B40").Value = 'Dai'
$ws.Range("C40").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E42").Value = '  -3.60%  '
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E43").Value = '  +4.96%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("E44").Value = '  -2.26%  '
$ws.Range("E45").Value = '  +1.34%  '
$ws.Range("E46").Value = '  -1.60%  '
$ws.Range("E47").Value = '  +17.44%  '
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("E49").Value = '  -3.20%  '
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("E51").Value = '  +2.49%  '

# Numeric-looking text values: must force text storage so Excel doesn't
# convert them to actual numbers (which would also normalize formatting,
# e.g. drop trailing zeros). We flip the cell to text format, write the
# value, then restore the default 'Normal' style so no stray formatting
# remains on the cell.
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '516.35'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '140.06'
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.604'
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.708'
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.167'
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0000320'
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '41.43'
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '10.24'
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '21.06'
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '13.92'
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '413.23'
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '3.44'
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '12.20'
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '13.91'
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '86.43'
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '3.98'
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.34'
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '35.22'
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '13.33'
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '675.83'
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.97'
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '66.27'
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.442'
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '39.22'
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.38'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.146'
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0472'
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.14'
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.85'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '3.45'
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.000281'
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.99'
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '3.28'
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '142.80'
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '8.69'
$c.Style = "Normal"
